# AutoCommit_14 июня 2024 г. 11:56:47_SibNout2023
# Update scores in the gradebook sheet (several students now have a 5 in
# columns that previously held lower marks / were empty), add a new
# homework column entry (лаб 1 follow-up note) and a note cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Батаев Вадим -------------------------------------------------
$ws.Range("C4").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5

# --- Row 5: Бесхлебный Даниэль -------------------------------------------
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 5

# J5 is a brand-new cell in this row; copy the number-column formatting
# (thick border style used by the rest of column J) from J4 before
# putting the value in, so the new cell keeps the sheet's styling.
$ws.Range("J4").Copy() | Out-Null
$ws.Range("J5").PasteSpecial(-4122) | Out-Null
$ws.Range("J5").Value = 5

# --- Row 15: Ильина Софья -------------------------------------------------
$ws.Range("H15").Value = 5
$ws.Range("I15").Value = 5

$ws.Range("J4").Copy() | Out-Null
$ws.Range("J15").PasteSpecial(-4122) | Out-Null
$ws.Range("J15").Value = 5

# --- Row 18: Коршунов Александр -------------------------------------------
# New J18 cell uses the alternate right-edge border style already used by
# J23 further down the column.
$ws.Range("J23").Copy() | Out-Null
$ws.Range("J18").PasteSpecial(-4122) | Out-Null
$ws.Range("J18").Value = 5

# --- Row 26: Полторабатько Кирилл -----------------------------------------
$ws.Range("F26").Value = 5
$ws.Range("G26").Value = 5
$ws.Range("I26").Value = 5

$ws.Range("J4").Copy() | Out-Null
$ws.Range("J26").PasteSpecial(-4122) | Out-Null
$ws.Range("J26").Value = 5

# Note cell for Полторабатько Кирилл ("очень ннада 4" == "really need a 4").
$ws.Range("N26").Value = "очень ннада 4"

$excel.CutCopyMode = 0

# --- View: scroll the frozen grid down a bit and leave the selection on
# the note cell we just added, matching where the editor ended up.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("C4").Select() | Out-Null
$win.FreezePanes = $true
$ws.Range("N28").Select() | Out-Null
